# Populate two more rows of username/password data on the active sheet
# (Sheet1), mirroring the existing header-less username/password layout
# in columns A:B, then move the selection the way the workbook shipped
# (cell D6 ends up selected after the edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ABC"
$ws.Range("B2").Value = "PASS123"
$ws.Range("A3").Value = "XYZ"
$ws.Range("B3").Value = "PASS456"

$ws.Range("D6").Select()
